# ============================================================
# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" (fund-holdings detail, same
#    layout as the existing 2021-Qx sheets) positioned right before
#    the "总计" (totals) summary sheet.
# 2. Prepend a new summary row for "2022-Q1" at the top of the "总计"
#    sheet's data table (shifting the existing rows down).
# ============================================================

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Step 1: create the "2022-Q1" worksheet, positioned before "总计"
# ------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q1Sheet = $wb.Worksheets.Add($totalSheet)
$q1Sheet.Name = "2022-Q1"

# NOTE: after Add()+rename, the $totalSheet variable captured above
# aliases the *same* anchor as $q1Sheet (both resolve to the freshly
# inserted sheet) - re-fetch a fresh, distinct reference to "总计" so
# the two sheets are no longer confused with one another.
$totalSheet = $wb.Worksheets.Item("总计")

# Use an existing quarter sheet as a formatting template: its header
# row (bold + bordered) and its first-column "index" style (bold +
# bordered) are copied over so the new sheet matches the others.
$template = $wb.Worksheets.Item("2021-Q4")

# Header row B1:H1 (bold, bordered, centered style)
$template.Range("B1:H1").Copy()
$q1Sheet.Range("B1:H1").PasteSpecial(-4122)

# Column-A "row index" style, stamped down for all 31 data rows
$template.Range("A2").Copy()
$q1Sheet.Range("A2:A32").PasteSpecial(-4122)

# Header labels
$q1Sheet.Cells.Item(1,2).Value = "基金代码"
$q1Sheet.Cells.Item(1,3).Value = "基金名称"
$q1Sheet.Cells.Item(1,4).Value = "基金规模"
$q1Sheet.Cells.Item(1,5).Value = "股票总仓位"
$q1Sheet.Cells.Item(1,6).Value = "仓位占比"
$q1Sheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$q1Sheet.Cells.Item(1,8).Value = "仓位排名"

# Fund-holdings detail rows: code, name, scale, total position,
# position ratio, market value held, position rank.
# Columns D/E/F/G hold numeric-looking values that are stored as TEXT
# in the source data, so they are entered with a leading apostrophe
# to force a text cell (matching "持有市值" etc. being plain numbers
# like "56.47" stored as strings, not numbers).
$q1Data = @(
    @("010347","农银汇理策略收益一年持有期混合","56.47","74.89","5.99","3.3826",4),
    @("100020","富国天益价值混合A","76.78","93.41","4.17","3.2017",10),
    @("270002","广发稳健增长混合A","178.96","41.88","1.74","3.1139",10),
    @("160926","大成创业板两年定期开放混合A","40.93","64.09","7.41","3.0329",1),
    @("660010","农银策略精选混合","31.49","75.60","6.73","2.1193",3),
    @("161040","富国创业板两年定期开放混合","35.36","83.63","4.31","1.5240",5),
    @("010815","农银汇理新兴消费股票","29.81","83.28","4.97","1.4816",5),
    @("000127","农银行业领先混合","12.28","75.75","6.49","0.7970",3),
    @("160529","博时创业板两年定期开放混合","7.92","82.61","8.87","0.7025",1),
    @("009798","大成创业板两年定期开放混合C","5.67","64.09","7.41","0.4201",1),
    @("008819","农银汇理策略趋势混合","6.17","76.46","6.67","0.4115",3),
    @("160143","南方创业板2年定期开放混合","8.23","84.71","4.67","0.3843",4),
    @("012260","广发睿明优质企业混合型证券投资基金A","12.09","61.24","2.96","0.3579",6),
    @("420003","天弘永定价值成长混合","6.55","81.42","3.55","0.2325",10),
    @("660012","农银消费主题混合A","6.36","65.93","3.34","0.2124",5),
    @("960033","农银消费主题混合H","6.36","65.93","3.34","0.2124",5),
    @("270041","广发消费品精选混合A","3.74","74.34","5.03","0.1881",4),
    @("164205","天弘文化新兴产业股票","4.45","82.58","4.21","0.1873",9),
    @("011078","诺德品质消费6个月持有期混合","4.49","93.14","3.91","0.1756",6),
    @("010457","广发睿鑫混合A","4.11","60.09","3.29","0.1352",4),
    @("002189","农银汇理国企改革灵活配置混合","1.15","54.15","7.66","0.0881",2),
    @("660003","农银平衡双利混合","3.75","69.47","1.94","0.0728",10),
    @("010808","达诚策略先锋混合A","0.73","84.86","7.01","0.0512",4),
    @("009326","广发稳健增长混合C","2.89","41.88","1.74","0.0503",10),
    @("011307","富国天益价值混合C","0.78","93.41","4.17","0.0325",10),
    @("010809","达诚策略先锋混合C","0.36","84.86","7.01","0.0252",4),
    @("012261","广发睿明优质企业混合型证券投资基金C","0.84","61.24","2.96","0.0249",6),
    @("010458","广发睿鑫混合C","0.64","60.09","3.29","0.0211",4),
    @("003308","中信建投睿利灵活配置混合A","0.08","93.35","9.46","0.0076",1),
    @("010022","广发消费品精选混合C","0.12","74.34","5.03","0.0060",4),
    @("004635","中信建投睿利灵活配置混合C","0.02","93.35","9.46","0.0019",1)
)

$r = 2
foreach ($row in $q1Data) {
    $q1Sheet.Cells.Item($r,1).Value = $r - 2
    $q1Sheet.Cells.Item($r,2).Value = "'" + $row[0]
    $q1Sheet.Cells.Item($r,3).Value = $row[1]
    $q1Sheet.Cells.Item($r,4).Value = "'" + $row[2]
    $q1Sheet.Cells.Item($r,5).Value = "'" + $row[3]
    $q1Sheet.Cells.Item($r,6).Value = "'" + $row[4]
    $q1Sheet.Cells.Item($r,7).Value = "'" + $row[5]
    $q1Sheet.Cells.Item($r,8).Value = $row[6]
    $r++
}

# The leading-apostrophe entries above force Excel to treat the
# (numeric-looking) values as text, but that also stamps a
# "quote-prefixed" style on those cells. Strip it back off by pasting
# in the (default/no-style) formatting from a guaranteed-blank region
# of this same new sheet, leaving the text values untouched.
$q1Sheet.Range("Z100:Z103").Copy()
$q1Sheet.Range("B2:B32").PasteSpecial(-4122)
$q1Sheet.Range("Z100:Z103").Copy()
$q1Sheet.Range("D2:G32").PasteSpecial(-4122)

# ------------------------------------------------------------------
# Step 2: prepend the "2022-Q1" summary row on the "总计" sheet
# ------------------------------------------------------------------
# Capture the existing 3 data rows (2021-Q4 / 2021-Q3 / 2021-Q2)
# before shifting them down by one row.
$oldDate = @("", "", "")
$oldCount = @(0, 0, 0)
$oldValue = @(0, 0, 0)
for ($i = 0; $i -lt 3; $i++) {
    $oldDate[$i] = $totalSheet.Cells.Item($i + 2, 2).Value2
    $oldCount[$i] = $totalSheet.Cells.Item($i + 2, 3).Value2
    $oldValue[$i] = $totalSheet.Cells.Item($i + 2, 4).Value2
}

# Extend the first-column index style down to the new last row (row 5)
$totalSheet.Range("A4").Copy()
$totalSheet.Range("A5").PasteSpecial(-4122)

# Write the shifted rows (old row 2 -> row 3, old row 3 -> row 4, old
# row 4 -> row 5), re-numbering the 0-based index column as we go.
for ($i = 2; $i -ge 0; $i--) {
    $destRow = $i + 3
    $totalSheet.Cells.Item($destRow, 1).Value = $i + 1
    $totalSheet.Cells.Item($destRow, 2).Value = $oldDate[$i]
    $totalSheet.Cells.Item($destRow, 3).Value = $oldCount[$i]
    $totalSheet.Cells.Item($destRow, 4).Value = $oldValue[$i]
}

# New "2022-Q1" row at the top of the data table
$totalSheet.Cells.Item(2,1).Value = 0
$totalSheet.Cells.Item(2,2).Value = "2022-Q1"
$totalSheet.Cells.Item(2,3).Value = 31
$totalSheet.Cells.Item(2,4).Value = 22.65

Write-Output "2022-Q1 sheet added; 总计 sheet updated"
